$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A97").Value = "2025/12/06 20:00"
$ws.Range("B97").Value = "-"
$ws.Range("C97").Value = "-"
$ws.Range("D97").Value = "-"
$ws.Range("E97").Value = "-"
$ws.Range("F97").Value = "-"
$ws.Range("G97").Value = "-"
